# Apply the "t-shirt size" column addition + Details text tweak + selection move
# to the WIP worksheet, matching the target commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: "t-shirt size" header + cycling S/M/L values for each data row ---
$ws.Range("F2").Value = "t-shirt size"

$sizes = @("S", "M", "L")
for ($row = 3; $row -le 21; $row++) {
    $ws.Cells.Item($row, 6).Value = $sizes[($row - 3) % 3]
}

# --- Column F width ---
# ColumnWidth is specified in character-width units; 15.1 round-trips to the
# stored OOXML width of exactly 16 (matching the target column width).
$ws.Columns.Item(6).ColumnWidth = 15.1

# --- Updated "Details" text for the Digital Display Phase 1 row (row 5) ---
$ws.Range("E5").Value = "50% approved, 50% in review, working on Phase 2, this is also test how much the text is flowing"

# --- Move the active selection to E10 ---
$ws.Range("E10").Select() | Out-Null
